$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(
    "17-09-2021",
    "18-09-2021",
    "19-09-2021",
    "20-09-2021",
    "21-09-2021",
    "22-09-2021",
    "23-09-2021",
    "24-09-2021",
    "25-09-2021",
    "26-09-2021",
    "27-09-2021",
    "28-09-2021",
    "29-09-2021",
    "30-09-2021"
)

$bValues = @(17537, 17537, 17537, 17537, 17537, 17537, 17537, 17537, 17537, 17537, 17537, 17537, 17537, 17537)
$cValues = @(1456, 1456, 1456, 1456, 1456, 1456, 1456, 1456, 1456, 1456, 1456, 737, 737, 737)
$dValues = @(521, 521, 521, 521, 521, 521, 521, 521, 521, 521, 521, 521, 521, 521)

$startRow = 261
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
    $ws.Cells.Item($row, 3).Value = $cValues[$i]
    $ws.Cells.Item($row, 4).Value = $dValues[$i]
}
